$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared string used by the "primary key" header column
# (USER_REFINED_PK -> SOURCE_TABLE_PRIMARY_KEY). The cell currently
# holding that text is X1; updating its value in place renames the
# shared-string table entry since it is the sole user of that string.
$ws.Range("X1").Value = "SOURCE_TABLE_PRIMARY_KEY"

# Move that column (now "SOURCE_TABLE_PRIMARY_KEY") from the end (X)
# to just left of the TARGET_* columns (I), shifting I:W -> J:X.
$ws.Columns("X:X").Cut()
$ws.Columns("I:I").Insert()

# The newly inserted column inherits the width of its left neighbour (H)
# automatically from a plain Insert, but Cut+Insert instead carries the
# cut column's own former width. Re-apply H's width onto the new column I
# so it matches the "insert column" formatting behaviour.
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth()

# Update the view: scroll the sheet so column C is left-most visible,
# and select I7 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("I7").Select() | Out-Null
